$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: mark new columns I:L as "spare" IOs (copy format from H5) ---
$ws.Range("I5").Value = "spare"
$ws.Range("J5").Value = "spare"
$ws.Range("K5").Value = "spare"
$ws.Range("L5").Value = "spare"
$ws.Range("H5").Copy()
$ws.Range("I5:L5").PasteSpecial(-4122)

# --- New row 32: Diagnostic_Header register ---
# (value entry order matches the shared-string insertion order of the target)
$ws.Range("B32").Value = "Diagnostic_Header"
$ws.Range("A32").Value = "0x001E"
$ws.Range("C32").Value = "RO"
$ws.Range("D32").Value = "Diagnostic Header IOs,Input:FPGA1-4+Teensy_FPGA_SP0-2 ,output: FPGA5-13"

# reg/function/R-W/description columns -> copy format from the row above (A31:D31)
$ws.Range("A31:D31").Copy()
$ws.Range("A32:D32").PasteSpecial(-4122)

# Input/output IO labels
$ws.Range("E32").Value = "FPGA1"
$ws.Range("F32").Value = "FPGA2"
$ws.Range("G32").Value = "FPGA3"
$ws.Range("H32").Value = "FPGA4"
$ws.Range("I32").Value = "Teensy_FPGA_SP0"
$ws.Range("J32").Value = "Teensy_FPGA_SP1"
$ws.Range("K32").Value = "Teensy_FPGA_SP2"
$ws.Range("L32").Value = "FPGA5"
$ws.Range("M32").Value = "FPGA6"
$ws.Range("N32").Value = "FPGA7"
$ws.Range("O32").Value = "FPGA8"
$ws.Range("P32").Value = "FPGA9"
$ws.Range("Q32").Value = "FPGA10"
$ws.Range("R32").Value = "FPGA11"
$ws.Range("S32").Value = "FPGA12"
$ws.Range("T32").Value = "FPGA13"

# E32:H32 and L32:T32 get a thin left/right border (no top/bottom, no fill).
# Build the border on a single seed cell, then format-paint it onto the rest
# so only one new border/style definition is minted.
$seed = $ws.Cells.Item(32, 5)
$seed.Borders.Item(7).Weight = 2
$seed.Borders.Item(10).Weight = 2

$seed.Copy()
$ws.Range("F32:H32").PasteSpecial(-4122)
$ws.Range("L32:T32").PasteSpecial(-4122)

$ws.Range("D32").Select()

$ws.PageSetup.Orientation = 1
